$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.3512740669568757
$ws.Range("C2").Value = 0.05555679239114397
$ws.Range("E2").Value = 0.4166729779671812
$ws.Range("F2").Value = 0.4443680307746121
$ws.Range("G2").Value = 0.002419048623702947
$ws.Range("K2").Value = 0.3305156662642332
$ws.Range("N2").Value = 1.265257978966236
$ws.Range("O2").Value = 2.156546062487251
$ws.Range("B3").Value = 0.3125479595428544
$ws.Range("C3").Value = 0.05062325806106571
$ws.Range("E3").Value = 0.3636028305573547
$ws.Range("F3").Value = 0.3878228170618172
$ws.Range("G3").Value = 0.002421454742344503
$ws.Range("K3").Value = 0.2895758957097314
$ws.Range("N3").Value = 1.280705205554181
$ws.Range("O3").Value = 2.170090714736901
$ws.Range("B4").Value = 0.2888069061714589
$ws.Range("C4").Value = 0.04757305651486377
$ws.Range("E4").Value = 0.3310981184881001
$ws.Range("F4").Value = 0.3531389305169483
$ws.Range("G4").Value = 0.002423008893010239
$ws.Range("K4").Value = 0.2644141422553901
$ws.Range("N4").Value = 1.290689865219008
$ws.Range("O4").Value = 2.179871185102044
$ws.Range("B5").Value = 0.2791418526768155
$ws.Range("C5").Value = 0.04632482103713187
$ws.Range("E5").Value = 0.3178709375501683
$ws.Range("F5").Value = 0.3390132514313251
$ws.Range("G5").Value = 0.002423661590890319
$ws.Range("K5").Value = 0.2541546515922164
$ws.Range("N5").Value = 1.294884386260689
$ws.Range("O5").Value = 2.184224201047954
$ws.Range("B6").Value = 0.2775375686617565
$ws.Range("C6").Value = 0.04611723592950057
$ws.Range("E6").Value = 0.3156756629712163
$ws.Range("F6").Value = 0.336668177824194
$ws.Range("G6").Value = 0.002423771142475449
$ws.Range("K6").Value = 0.2524507254285879
$ws.Range("N6").Value = 1.2955884736696
$ws.Range("O6").Value = 2.184969184209294
$ws.Range("B7").Value = 0.2886765204945334
$ws.Range("C7").Value = 0.04755624361403932
$ws.Range("E7").Value = 0.3309196581478204
$ws.Range("F7").Value = 0.3529483938344953
$ws.Range("G7").Value = 0.00242301761709659
$ws.Range("K7").Value = 0.2642758025434091
$ws.Range("N7").Value = 1.290745925130842
$ws.Range("O7").Value = 2.179928404818625
$ws.Range("B8").Value = 0.337913825893736
$ws.Range("C8").Value = 0.05386007997141462
$ws.Range("E8").Value = 0.3983567442039799
$ws.Range("F8").Value = 0.4248636149813478
$ws.Range("G8").Value = 0.002419862353590241
$ws.Range("K8").Value = 0.3164049033385368
$ws.Range("N8").Value = 1.270480385833268
$ws.Range("O8").Value = 2.160911956837694
$ws.Range("B9").Value = 0.4347514544876674
$ws.Range("C9").Value = 0.06605481153697212
$ws.Range("E9").Value = 0.531314160805465
$ws.Range("F9").Value = 0.5661985755041457
$ws.Range("G9").Value = 0.00241428135420324
$ws.Range("K9").Value = 0.418426599755179
$ws.Range("N9").Value = 1.234708184528206
$ws.Range("O9").Value = 2.135271886523782
$ws.Range("B10").Value = 0.5060651959437905
$ws.Range("C10").Value = 0.07491259671296291
$ws.Range("E10").Value = 0.6295445602724641
$ws.Range("F10").Value = 0.6702781546542269
$ws.Range("G10").Value = 0.002410546795283614
$ws.Range("K10").Value = 0.4932542622877634
$ws.Range("N10").Value = 1.210846508765403
$ws.Range("O10").Value = 2.123587643167411
$ws.Range("B11").Value = 0.5385432330129163
$ws.Range("C11").Value = 0.07892024778979589
$ws.Range("E11").Value = 0.6743754978692493
$ws.Range("F11").Value = 0.7176906081379002
$ws.Range("G11").Value = 0.002408926443263616
$ws.Range("K11").Value = 0.5272672814735984
$ws.Range("N11").Value = 1.200516957289967
$ws.Range("O11").Value = 2.119836653585139
$ws.Range("B12").Value = 0.550846924136664
$ws.Range("C12").Value = 0.08043469491450139
$ws.Range("E12").Value = 0.6913745423223219
$ws.Range("F12").Value = 0.7356546913071611
$ws.Range("G12").Value = 0.002408324085590126
$ws.Range("K12").Value = 0.5401431266747636
$ws.Range("N12").Value = 1.196681011236315
$ws.Range("O12").Value = 2.118642072982624
$ws.Range("B13").Value = 0.5481968891654958
$ws.Range("C13").Value = 0.08010867277124589
$ws.Range("E13").Value = 0.6877124653494349
$ws.Range("F13").Value = 0.7317853510981394
$ws.Range("G13").Value = 0.002408453315166863
$ws.Range("K13").Value = 0.5373702728924741
$ws.Range("N13").Value = 1.197503785491719
$ws.Range("O13").Value = 2.118889288653378
$ws.Range("B14").Value = 0.5395553677010696
$ws.Range("C14").Value = 0.07904490582788526
$ws.Range("E14").Value = 0.6757735583504143
$ws.Range("F14").Value = 0.7191683204515869
$ws.Range("G14").Value = 0.00240887666194769
$ws.Range("K14").Value = 0.5283266700704132
$ws.Range("N14").Value = 1.200199855349183
$ws.Range("O14").Value = 2.119733843403566
$ws.Range("B15").Value = 0.5342628218414802
$ws.Range("C15").Value = 0.07839290505582142
$ws.Range("E15").Value = 0.66846361718423
$ws.Range("F15").Value = 0.7114413442032514
$ws.Range("G15").Value = 0.002409137435517041
$ws.Range("K15").Value = 0.522786652039855
$ws.Range("N15").Value = 1.20186112843513
$ws.Range("O15").Value = 2.120280595209437
$ws.Range("B16").Value = 0.5039433942536675
$ws.Range("C16").Value = 0.07465024616676885
$ws.Range("E16").Value = 0.6266178184512796
$ws.Range("F16").Value = 0.6671810134426437
$ws.Range("G16").Value = 0.002410654264390136
$ws.Range("K16").Value = 0.4910308697572248
$ws.Range("N16").Value = 1.211532147088935
$ws.Range("O16").Value = 2.123864332577085
$ws.Range("B17").Value = 0.4853526382827908
$ws.Range("C17").Value = 0.07234864222481008
$ws.Range("E17").Value = 0.6009851798606292
$ws.Range("F17").Value = 0.6400460337125793
$ws.Range("G17").Value = 0.002411604859800914
$ws.Range("K17").Value = 0.4715427047674439
$ws.Range("N17").Value = 1.217599602443183
$ws.Range("O17").Value = 2.126464117527576
$ws.Range("B18").Value = 0.4746632370008115
$ws.Range("C18").Value = 0.07102277028812409
$ws.Range("E18").Value = 0.5862555551440636
$ws.Range("F18").Value = 0.6244449056556647
$ws.Range("G18").Value = 0.002412159010867295
$ws.Range("K18").Value = 0.4603311574923907
$ws.Range("N18").Value = 1.221138874462255
$ws.Range("O18").Value = 2.128106623562246
$ws.Range("B19").Value = 0.4710446033867299
$ws.Range("C19").Value = 0.07057350158348186
$ws.Range("E19").Value = 0.581270642188926
$ws.Range("F19").Value = 0.619163680173358
$ws.Range("G19").Value = 0.002412347908526873
$ws.Range("K19").Value = 0.4565347064253444
$ws.Range("N19").Value = 1.22234569815862
$ws.Range("O19").Value = 2.128687998512959
$ws.Range("B20").Value = 0.4873312953915843
$ws.Range("C20").Value = 0.07259386429261383
$ws.Range("E20").Value = 0.6037123996433706
$ws.Range("F20").Value = 0.642933953830422
$ws.Range("G20").Value = 0.002411502902349518
$ws.Range("K20").Value = 0.473617511343889
$ws.Range("N20").Value = 1.216948594690258
$ws.Range("O20").Value = 2.126172127022187
$ws.Range("B21").Value = 0.5420934602291538
$ws.Range("C21").Value = 0.07935744604309036
$ws.Range("E21").Value = 0.6792796799325771
$ws.Range("F21").Value = 0.7228739723491628
$ws.Range("G21").Value = 0.002408752010419961
$ws.Range("K21").Value = 0.5309831107023228
$ws.Range("N21").Value = 1.19940590084995
$ws.Range("O21").Value = 2.119479640804315
$ws.Range("B22").Value = 0.5779124625587428
$ws.Range("C22").Value = 0.08375938630277346
$ws.Range("E22").Value = 0.7287994196467196
$ws.Range("F22").Value = 0.7751780083420101
$ws.Range("G22").Value = 0.002407019601227142
$ws.Range("K22").Value = 0.5684505171652461
$ws.Range("N22").Value = 1.18838162088873
$ws.Range("O22").Value = 2.11642249234103
$ws.Range("B23").Value = 0.5587926754835735
$ws.Range("C23").Value = 0.08141168556092282
$ws.Range("E23").Value = 0.7023571690165511
$ws.Range("F23").Value = 0.7472568307830727
$ws.Range("G23").Value = 0.002407938249679989
$ws.Range("K23").Value = 0.5484557957257152
$ws.Range("N23").Value = 1.1942251111339
$ws.Range("O23").Value = 2.117933364932838
$ws.Range("B24").Value = 0.4864367486226229
$ws.Range("C24").Value = 0.07248300764634052
$ws.Range("E24").Value = 0.6024794021762716
$ws.Range("F24").Value = 0.6416283278902171
$ws.Range("G24").Value = 0.002411548973494201
$ws.Range("K24").Value = 0.4726795147581697
$ws.Range("N24").Value = 1.21724275639772
$ws.Range("O24").Value = 2.126303675419393
$ws.Range("B25").Value = 0.408524558834614
$ws.Range("C25").Value = 0.06277368414033901
$ws.Range("E25").Value = 0.4952577799044775
$ws.Range("F25").Value = 0.5279251897347166
$ws.Range("G25").Value = 0.002415726642553103
$ws.Range("K25").Value = 0.3908491882121155
$ws.Range("N25").Value = 1.243960585796582
$ws.Range("O25").Value = 2.140955653543358
